# Append: 2026-01-13 06:31 JST
# Replaces the 13-row job listing with a fresh 7-row listing scraped at
# 2026-01-13 06:31:06, narrowing column B and dropping the now-stale rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$timestamp = "2026-01-13 06:31:06"

# --- New row data (rows 2-8) --------------------------------------------
# title, category, price, deadline, url, score, skills
$rows = @(
    @("【募集】Python / Docker 日次データ スクレイピングシステム構築", "システム開発", "50,000 円 ~ 100,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5469627", 248, "🔥Python ◆スクレイピング"),
    @("【急募】FXツール開発のプロフェッショナルを探しています!", "システム開発", "20,000 円 ~ 50,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5470011", 123, "◆ツール,開発"),
    @("自社システムの開発・保守エンジニア募集★☆カメラ面談あり", "システム開発", "5,000 円 ~ 10,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5469878", 75, "◆開発"),
    @("初回 東京スポーツ施設の空きテニスコートを自動予約してくれるシステムの開発", "システム開発", "10,000 円 ~ 20,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5469843", 75, "◆開発"),
    @("GoogleCloudを利用したアジャイル開発共通基盤のSREエンジニアの募集", "システム開発", "500,000 円 ~ 1,000,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5457458", 75, "◆開発"),
    @("【緊急】AWS上の稼働中Webサイトを最新GitHubに再構築", "システム開発", "5,000 円 ~ 10,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5469840", 30, "◇サイト"),
    @("n8n 初期構築・セットアップ(小規模/検証用途)", "システム開発", "10,000 円 ~ 20,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5469826", 10, "")
)

# --- Drop every existing hyperlink up front ------------------------------
# (the runtime's Hyperlinks.Delete() clears the whole sheet collection no
# matter which range it's invoked on, so do this exactly once before
# rebuilding the F-column links we actually want to keep)
$ws.Range("A1").Hyperlinks.Delete()

# --- Remove the now-stale rows 9:14 --------------------------------------
$ws.Rows("9:14").Delete()

# --- Narrow column B ------------------------------------------------------
# ColumnWidth uses Excel's "characters" unit, which round-trips to the
# stored OOXML width via +5/6; back that off so the saved width is exactly
# 41 (matching the authored diff).
$ws.Columns.Item(2).ColumnWidth = 40.166666666666664

# --- Write rows 2-8 and re-link column F ---------------------------------
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $data = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $timestamp
    $ws.Cells.Item($r, 2).Value = $data[0]
    $ws.Cells.Item($r, 3).Value = $data[1]
    $ws.Cells.Item($r, 4).Value = $data[2]
    $ws.Cells.Item($r, 5).Value = $data[3]
    $ws.Cells.Item($r, 6).Value = $data[4]
    $ws.Cells.Item($r, 7).Value = $data[5]
    $ws.Cells.Item($r, 8).Value = $data[6]

    $linkCell = $ws.Cells.Item($r, 6)
    $ws.Hyperlinks.Add($linkCell, $data[4])
    # Hyperlinks.Add silently re-styles the cell with a duplicate xf; force
    # it back onto the canonical "Hyperlink" style so it matches the
    # original styling (cellXfs index 1).
    $linkCell.Style = "Hyperlink"
}

# Row 8's skill-summary column has no value in the new data -- make sure it
# is truly empty (no cell entry at all), matching the authored sheet.
$ws.Cells.Item(8, 8).ClearContents()
